{"js": "// Update the answers in the \"two-digit number divided by one-digit number\"\n// practice table. Each data row of the table (rows 0, 4, 8, 12, 16 - the\n// rows in between are spacer/blank rows) holds 5 division problems, one per\n// column. We overwrite the text of each of those 25 cells with its new\n// value while leaving every other part of the document (formatting,\n// paragraph/run properties, blank rows, etc.) untouched.\nconst updates = [\n  { row: 0, col: 0, before: \"45\u00f79=5, 0\", after: \"33\u00f72=16, 1\" },\n  { row: 0, col: 1, before: \"53\u00f74=13, 1\", after: \"42\u00f77=6, 0\" },\n  { row: 0, col: 2, before: \"63\u00f76=10, 3\", after: \"10\u00f75=2, 0\" },\n  { row: 0, col: 3, before: \"62\u00f79=6, 8\", after: \"89\u00f74=22, 1\" },\n  { row: 0, col: 4, before: \"46\u00f78=5, 6\", after: \"62\u00f77=8, 6\" },\n  { row: 4, col: 0, before: \"50\u00f77=7, 1\", after: \"66\u00f76=11, 0\" },\n  { row: 4, col: 1, before: \"89\u00f75=17, 4\", after: \"92\u00f79=10, 2\" },\n  { row: 4, col: 2, before: \"17\u00f74=4, 1\", after: \"57\u00f75=11, 2\" },\n  { row: 4, col: 3, before: \"74\u00f72=37, 0\", after: \"50\u00f73=16, 2\" },\n  { row: 4, col: 4, before: \"48\u00f75=9, 3\", after: \"15\u00f78=1, 7\" },\n  { row: 8, col: 0, before: \"19\u00f75=3, 4\", after: \"80\u00f75=16, 0\" },\n  { row: 8, col: 1, before: \"62\u00f72=31, 0\", after: \"96\u00f77=13, 5\" },\n  { row: 8, col: 2, before: \"65\u00f72=32, 1\", after: \"82\u00f75=16, 2\" },\n  { row: 8, col: 3, before: \"97\u00f79=10, 7\", after: \"87\u00f77=12, 3\" },\n  { row: 8, col: 4, before: \"80\u00f75=16, 0\", after: \"48\u00f72=24, 0\" },\n  { row: 12, col: 0, before: \"43\u00f74=10, 3\", after: \"11\u00f72=5, 1\" },\n  { row: 12, col: 1, before: \"56\u00f79=6, 2\", after: \"59\u00f72=29, 1\" },\n  { row: 12, col: 2, before: \"39\u00f79=4, 3\", after: \"57\u00f77=8, 1\" },\n  { row: 12, col: 3, before: \"25\u00f72=12, 1\", after: \"92\u00f79=10, 2\" },\n  { row: 12, col: 4, before: \"82\u00f73=27, 1\", after: \"98\u00f75=19, 3\" },\n  { row: 16, col: 0, before: \"94\u00f74=23, 2\", after: \"50\u00f76=8, 2\" },\n  { row: 16, col: 1, before: \"36\u00f76=6, 0\", after: \"60\u00f72=30, 0\" },\n  { row: 16, col: 2, before: \"10\u00f77=1, 3\", after: \"74\u00f74=18, 2\" },\n  { row: 16, col: 3, before: \"15\u00f74=3, 3\", after: \"28\u00f74=7, 0\" },\n  { row: 16, col: 4, before: \"44\u00f77=6, 2\", after: \"56\u00f72=28, 0\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Grab the first paragraph range of every target cell up front (one sync),\n// then mutate each range's text (another sync) so formatting (rFonts/sz,\n// paragraph alignment, etc.) carried by the existing run is preserved -\n// only the <w:t> text itself changes, exactly like the source diff.\nconst ranges = [];\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items/text\");\n  ranges.push({ u, paragraphs });\n}\nawait context.sync();\n\nfor (const { u, paragraphs } of ranges) {\n  const para = paragraphs.items[0];\n  // Sanity-check we are editing the expected cell before overwriting it.\n  if (para.text !== u.before) {\n    throw new Error(\n      `Unexpected cell text at row ${u.row}, col ${u.col}: ` +\n        `expected \"${u.before}\" but found \"${para.text}\"`\n    );\n  }\n  const range = para.getRange();\n  range.insertText(u.after, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the answers in the \"two-digit number divided by one-digit number\"\n# practice table. Each data row of the table (1-indexed rows 1, 5, 9, 13, 17\n# - the rows in between are spacer/blank rows) holds 5 division problems,\n# one per column. We overwrite the text of each of those 25 cells with its\n# new value while leaving every other part of the document (formatting,\n# paragraph/run properties, blank rows, etc.) untouched.\n\n$updates = @(\n  @{ Row = 1;  Col = 1; Before = \"45\u00f79=5, 0\";   After = \"33\u00f72=16, 1\" },\n  @{ Row = 1;  Col = 2; Before = \"53\u00f74=13, 1\";  After = \"42\u00f77=6, 0\" },\n  @{ Row = 1;  Col = 3; Before = \"63\u00f76=10, 3\";  After = \"10\u00f75=2, 0\" },\n  @{ Row = 1;  Col = 4; Before = \"62\u00f79=6, 8\";   After = \"89\u00f74=22, 1\" },\n  @{ Row = 1;  Col = 5; Before = \"46\u00f78=5, 6\";   After = \"62\u00f77=8, 6\" },\n\n  @{ Row = 5;  Col = 1; Before = \"50\u00f77=7, 1\";   After = \"66\u00f76=11, 0\" },\n  @{ Row = 5;  Col = 2; Before = \"89\u00f75=17, 4\";  After = \"92\u00f79=10, 2\" },\n  @{ Row = 5;  Col = 3; Before = \"17\u00f74=4, 1\";   After = \"57\u00f75=11, 2\" },\n  @{ Row = 5;  Col = 4; Before = \"74\u00f72=37, 0\";  After = \"50\u00f73=16, 2\" },\n  @{ Row = 5;  Col = 5; Before = \"48\u00f75=9, 3\";   After = \"15\u00f78=1, 7\" },\n\n  @{ Row = 9;  Col = 1; Before = \"19\u00f75=3, 4\";   After = \"80\u00f75=16, 0\" },\n  @{ Row = 9;  Col = 2; Before = \"62\u00f72=31, 0\";  After = \"96\u00f77=13, 5\" },\n  @{ Row = 9;  Col = 3; Before = \"65\u00f72=32, 1\";  After = \"82\u00f75=16, 2\" },\n  @{ Row = 9;  Col = 4; Before = \"97\u00f79=10, 7\";  After = \"87\u00f77=12, 3\" },\n  @{ Row = 9;  Col = 5; Before = \"80\u00f75=16, 0\";  After = \"48\u00f72=24, 0\" },\n\n  @{ Row = 13; Col = 1; Before = \"43\u00f74=10, 3\";  After = \"11\u00f72=5, 1\" },\n  @{ Row = 13; Col = 2; Before = \"56\u00f79=6, 2\";   After = \"59\u00f72=29, 1\" },\n  @{ Row = 13; Col = 3; Before = \"39\u00f79=4, 3\";   After = \"57\u00f77=8, 1\" },\n  @{ Row = 13; Col = 4; Before = \"25\u00f72=12, 1\";  After = \"92\u00f79=10, 2\" },\n  @{ Row = 13; Col = 5; Before = \"82\u00f73=27, 1\";  After = \"98\u00f75=19, 3\" },\n\n  @{ Row = 17; Col = 1; Before = \"94\u00f74=23, 2\";  After = \"50\u00f76=8, 2\" },\n  @{ Row = 17; Col = 2; Before = \"36\u00f76=6, 0\";   After = \"60\u00f72=30, 0\" },\n  @{ Row = 17; Col = 3; Before = \"10\u00f77=1, 3\";   After = \"74\u00f74=18, 2\" },\n  @{ Row = 17; Col = 4; Before = \"15\u00f74=3, 3\";   After = \"28\u00f74=7, 0\" },\n  @{ Row = 17; Col = 5; Before = \"44\u00f77=6, 2\";   After = \"56\u00f72=28, 0\" }\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nforeach ($u in $updates) {\n  $cell = $t.Cell($u.Row, $u.Col)\n  $range = $cell.Range\n  # A cell range's text ends with the cell-end mark (CR + cell marker);\n  # strip it off before comparing so the sanity check lines up with the\n  # plain text value we expect to find.\n  $current = $range.Text.TrimEnd([char]13, [char]7)\n  if ($current -ne $u.Before) {\n    throw \"Unexpected cell text at row $($u.Row), col $($u.Col): expected '$($u.Before)' but found '$current'\"\n  }\n  $range.Text = $u.After\n}\n"}
